$wb = $excel.ActiveWorkbook

# --- Sheet "2025-1": add the two new fleet entries in rows 18 and 19 ---
$ws1 = $wb.Worksheets.Item("2025-1")

$ws1.Range("A18").Value = "2025-1"
$ws1.Range("B18").Value = "EP TASA 413 ARBOLADURA"
$ws1.Range("C18").Value = "Embarcación Pesquera"
$ws1.Range("D18").Value = "GP/83"
$ws1.Range("E18").Value = "GP/83-125"

$ws1.Range("A19").Value = "2025-1"
$ws1.Range("D19").Value = "GP/84"
$ws1.Range("B19").Value = "EP TASA 424 ARBOLADURA"
$ws1.Range("C19").Value = "Embarcación Pesquera"
$ws1.Range("E19").Value = "GP/84-125"

# Widen column E a bit to fit the new content
$ws1.Columns.Item(5).ColumnWidth = 13.3

# Move the active selection to F18 (matches where the editor left off)
$ws1.Range("F18").Select()

# --- Sheet "2024-2": scroll the view down without disturbing the existing selection ---
$ws2 = $wb.Worksheets.Item("2024-2")
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 9

# Restore "2025-1" as the active/selected sheet (it was the tab shown before saving)
$ws1.Activate()
$ws1.Range("F18").Select()
